$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should match the formatting of the existing header row (A1:AC1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (Wins / Losses / Ties) repeated for every player row
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 90   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 72   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
